# Update auto scs (lamda_1, column B) and time in ms (lamda_2, column C)
# for every data row, update the dic_nbre_clients_poisson_2_keys (D) and
# dic_nbre_clients_prob_poisson_2_values (E) columns, and append a new
# data row (row 55) with the "auto capacity" data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the new row (row 55) before filling in data, copying the
#        formatting (bold / centered / bordered) from the last existing
#        data row (row 54) so the new row matches the existing ones.
$ws.Cells.Item(54, 1).Copy()
$ws.Cells.Item(55, 1).PasteSpecial(-4122)

# --- 2. New D (keys) / E (values) data for rows 2..55
$dVals = @(0,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,48,49,50,51,53,54,55,56,59)
$eVals = @(0.133,0.001,0.004,0.008,0.026,0.037,0.044,0.039,0.041,0.035,0.03,0.032,0.021,0.037,0.038,0.041,0.045,0.041,0.034,0.023,0.03,0.018,0.02,0.018,0.021,0.012,0.02,0.017,0.014,0.009000000000000001,0.01,0.014,0.009000000000000001,0.008,0.007,0.005,0.009000000000000001,0.002,0.006,0.005,0.004,0.005,0.005,0.003,0.002,0.003,0.002,0.002,0.001,0.001,0.001,0.001,0.001,0.001)

$firstRow = 2
$lastRow = 55

# --- 3. Update lamda_1 (B) and lamda_2 (C) for every data row, and the
#        D / E columns from the arrays above.
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 2).Value = 33.94444444444444
    $ws.Cells.Item($r, 3).Value = 1.95
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
    $ws.Cells.Item($r, 5).Value = $eVals[$i]
}

# --- 4. Set the key value for the newly added row (A55 = 53, continuing
#        the existing A column sequence).
$ws.Cells.Item($lastRow, 1).Value = $lastRow - 2
